# Refresh the "cryptos" price/volume snapshot (and fix the row order for the
# OKB / Stacks / Bittensor trio) to match the latest GitHub Actions pull.
#
# Price values like "601.80" look numeric to Excel's automatic type
# detection, so for those cells we briefly force a Text number format
# before writing the string and then restore the cell's original
# (default/"Normal") style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.119.91"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.810.22"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("D7").Value = "3.806.83"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "4.448.35"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "3.840.25"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "69.225.54"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("D33").Value = "3.959.30"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "3.756.64"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("E37").Value = "  +5.13%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "429.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "2.836.10"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0351"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
